# Daily attendance processing - 2025-10-03 06:24:43
# Applies the data corrections captured in the commit diff to the
# "Session Analysis Results" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row-level status changes (row fill/style reflects the Status column)
#    Style 2 = Pending (light yellow), Style 4 = Not Recorded (pink),
#    Style 6 = Recorded (green). We copy the *format* from an existing
#    row that already carries the desired style so the workbook keeps
#    reusing the same style definitions instead of inventing new ones.
# ---------------------------------------------------------------------

# Row 57: Not Recorded -> Pending
$ws.Range("A2:I2").Copy()
$ws.Range("A57:I57").PasteSpecial(-4122)

# Row 77: Pending -> Not Recorded
$ws.Range("A3:I3").Copy()
$ws.Range("A77:I77").PasteSpecial(-4122)

# Row 83: Recorded -> Not Recorded
$ws.Range("A3:I3").Copy()
$ws.Range("A83:I83").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Content changes tied to the three re-styled rows
# ---------------------------------------------------------------------

# Row 57 - date moved a session from 01/10/2025 to 10/11/2025; now pending
$ws.Range("E57").NumberFormat = "@"
$ws.Range("E57").Value = "10/11/2025"
$ws.Range("I57").Value = "Pending"

# Row 77 - date corrected from 07/10/2025 to 10/07/2025; now not recorded
$ws.Range("E77").NumberFormat = "@"
$ws.Range("E77").Value = "10/07/2025"
$ws.Range("I77").Value = "Not Recorded"

# Row 83 - date corrected to 10/01/2025; recording was reverted/cleared
$ws.Range("E83").NumberFormat = "@"
$ws.Range("E83").Value = "10/01/2025"
$ws.Range("G83").Value = ""
$ws.Range("H83").Value = "0/154"
$ws.Range("I83").Value = "Not Recorded"

# Row 78 - date corrected from 12/10/2025 to 10/12/2025
$ws.Range("E78").NumberFormat = "@"
$ws.Range("E78").Value = "10/12/2025"

# ---------------------------------------------------------------------
# 3) Class Statistics summary numbers (columns K:L)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 39
$ws.Range("L7").Value = 10

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "25.5%"

$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "47.9%"

# ---------------------------------------------------------------------
# 4) Per-group breakdown table (columns K:S)
# ---------------------------------------------------------------------
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 11

$ws.Range("O19").Value = 4
$ws.Range("P19").Value = 2
$ws.Range("Q19").Value = 11

$ws.Range("R19").NumberFormat = "@"
$ws.Range("R19").Value = "23.5%"

$ws.Range("S19").NumberFormat = "@"
$ws.Range("S19").Value = "56.2%"

# ---------------------------------------------------------------------
# 5) "Recorded By" email list re-orderings
# ---------------------------------------------------------------------
$ws.Range("G17").Value = "Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G34").Value = "Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

$ws.Range("G45").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, System, backup@backdoor.com"
$ws.Range("G62").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, System, backup@backdoor.com"

$ws.Range("G51").Value = "Salma.hassan@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G68").Value = "Salma.hassan@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"

$ws.Range("G72").Value = "Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

$ws.Range("G76").Value = "mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"

$ws.Range("G85").Value = "neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G102").Value = "neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"

$ws.Range("G98").Value = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"

$ws.Range("G99").Value = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G149").Value = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg"

$ws.Range("G116").Value = "nourhan.mostafa@med.asu.edu.eg, enas.omran@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G133").Value = "nourhan.mostafa@med.asu.edu.eg, enas.omran@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"

$ws.Range("G119").Value = "neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G136").Value = "neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"

$ws.Range("G150").Value = "Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
